# BOT; UPDATE DATA
# Insert the newest daily data point (2020-04-16 / serial 43937) as a new
# row 82 on "相談件数", pushing the existing footnote row down to row 83,
# then refresh the print area / visible selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Insert a new row before the footnote row (currently row 82) and fill it
# with the latest day's figures. EntireRow insert shifts the footnote row
# (and its styles) down to row 83 automatically.
$ws.Rows.Item(82).Insert()

$ws.Range("A82").Value = 43937
$ws.Range("B82").Value = 716
$ws.Range("C82").Value = 24208
$ws.Range("D82").Value = 178
$ws.Range("E82").Value = 5548

# Grow the print area by one row so the new data row is included.
foreach ($n in $wb.Names) {
    $n.RefersTo = "=相談件数!`$A`$1:`$E`$87"
}

# Bring the newly-added row into view / update the remembered selection.
[void]$ws.Range("E84").Select()
